# Refresh scrape timestamps (scrapedAt / lastSeenAt) for the Bangalore events
# sheet after re-running the scraper. Event data (columns A-G) is unchanged;
# only the bookkeeping timestamp columns H and I move forward to the new run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; ScrapedAt = "2026-02-03T14:49:20.046Z"; LastSeenAt = "2026-02-03T14:49:20.075Z" }
    @{ Row = 3; ScrapedAt = "2026-02-03T14:49:20.047Z"; LastSeenAt = "2026-02-03T14:49:20.082Z" }
    @{ Row = 4; ScrapedAt = "2026-02-03T14:49:20.047Z"; LastSeenAt = "2026-02-03T14:49:20.084Z" }
    @{ Row = 5; ScrapedAt = "2026-02-03T14:49:20.047Z"; LastSeenAt = "2026-02-03T14:49:20.087Z" }
    @{ Row = 6; ScrapedAt = "2026-02-03T14:49:20.047Z"; LastSeenAt = "2026-02-03T14:49:20.095Z" }
    @{ Row = 7; ScrapedAt = "2026-02-03T14:49:20.048Z"; LastSeenAt = "2026-02-03T14:49:20.097Z" }
    @{ Row = 8; ScrapedAt = "2026-02-03T14:49:20.048Z"; LastSeenAt = "2026-02-03T14:49:20.101Z" }
    @{ Row = 9; ScrapedAt = "2026-02-03T14:49:20.048Z"; LastSeenAt = "2026-02-03T14:49:20.104Z" }
    @{ Row = 10; ScrapedAt = "2026-02-03T14:49:20.048Z"; LastSeenAt = "2026-02-03T14:49:20.107Z" }
    @{ Row = 11; ScrapedAt = "2026-02-03T14:49:20.049Z"; LastSeenAt = "2026-02-03T14:49:20.111Z" }
    @{ Row = 12; ScrapedAt = "2026-02-03T14:49:20.049Z"; LastSeenAt = "2026-02-03T14:49:20.118Z" }
    @{ Row = 13; ScrapedAt = "2026-02-03T14:49:20.049Z"; LastSeenAt = "2026-02-03T14:49:20.122Z" }
    @{ Row = 14; ScrapedAt = "2026-02-03T14:49:20.049Z"; LastSeenAt = "2026-02-03T14:49:20.124Z" }
    @{ Row = 15; ScrapedAt = "2026-02-03T14:49:20.049Z"; LastSeenAt = "2026-02-03T14:49:20.127Z" }
    @{ Row = 16; ScrapedAt = "2026-02-03T14:49:20.049Z"; LastSeenAt = "2026-02-03T14:49:20.131Z" }
    @{ Row = 17; ScrapedAt = "2026-02-03T14:49:20.049Z"; LastSeenAt = "2026-02-03T14:49:20.133Z" }
    @{ Row = 18; ScrapedAt = "2026-02-03T14:49:20.049Z"; LastSeenAt = "2026-02-03T14:49:20.136Z" }
    @{ Row = 19; ScrapedAt = "2026-02-03T14:49:20.049Z"; LastSeenAt = "2026-02-03T14:49:20.139Z" }
    @{ Row = 20; ScrapedAt = "2026-02-03T14:49:20.050Z"; LastSeenAt = "2026-02-03T14:49:20.142Z" }
    @{ Row = 21; ScrapedAt = "2026-02-03T14:49:20.050Z"; LastSeenAt = "2026-02-03T14:49:20.145Z" }
    @{ Row = 22; ScrapedAt = "2026-02-03T14:49:20.050Z"; LastSeenAt = "2026-02-03T14:49:20.147Z" }
    @{ Row = 23; ScrapedAt = "2026-02-03T14:49:20.050Z"; LastSeenAt = "2026-02-03T14:49:20.151Z" }
    @{ Row = 24; ScrapedAt = "2026-02-03T14:49:20.050Z"; LastSeenAt = "2026-02-03T14:49:20.158Z" }
    @{ Row = 25; ScrapedAt = "2026-02-03T14:49:20.050Z"; LastSeenAt = "2026-02-03T14:49:20.160Z" }
    @{ Row = 26; ScrapedAt = "2026-02-03T14:49:20.050Z"; LastSeenAt = "2026-02-03T14:49:20.162Z" }
    @{ Row = 27; ScrapedAt = "2026-02-03T14:49:20.050Z"; LastSeenAt = "2026-02-03T14:49:20.165Z" }
    @{ Row = 28; ScrapedAt = "2026-02-03T14:49:20.051Z"; LastSeenAt = "2026-02-03T14:49:20.167Z" }
    @{ Row = 29; ScrapedAt = "2026-02-03T14:49:20.051Z"; LastSeenAt = "2026-02-03T14:49:20.170Z" }
    @{ Row = 30; ScrapedAt = "2026-02-03T14:49:20.051Z"; LastSeenAt = "2026-02-03T14:49:20.173Z" }
    @{ Row = 31; ScrapedAt = "2026-02-03T14:49:20.051Z"; LastSeenAt = "2026-02-03T14:49:20.175Z" }
    @{ Row = 32; ScrapedAt = "2026-02-03T14:49:20.051Z"; LastSeenAt = "2026-02-03T14:49:20.178Z" }
    @{ Row = 33; ScrapedAt = "2026-02-03T14:49:20.052Z"; LastSeenAt = "2026-02-03T14:49:20.180Z" }
    @{ Row = 34; ScrapedAt = "2026-02-03T14:49:20.052Z"; LastSeenAt = "2026-02-03T14:49:20.182Z" }
    @{ Row = 35; ScrapedAt = "2026-02-03T14:49:20.052Z"; LastSeenAt = "2026-02-03T14:49:20.185Z" }
    @{ Row = 36; ScrapedAt = "2026-02-03T14:49:20.052Z"; LastSeenAt = "2026-02-03T14:49:20.187Z" }
    @{ Row = 37; ScrapedAt = "2026-02-03T14:49:20.052Z"; LastSeenAt = "2026-02-03T14:49:20.190Z" }
    @{ Row = 38; ScrapedAt = "2026-02-03T14:49:20.052Z"; LastSeenAt = "2026-02-03T14:49:20.191Z" }
    @{ Row = 39; ScrapedAt = "2026-02-03T14:49:20.052Z"; LastSeenAt = "2026-02-03T14:49:20.194Z" }
    @{ Row = 40; ScrapedAt = "2026-02-03T14:49:20.053Z"; LastSeenAt = "2026-02-03T14:49:20.196Z" }
    @{ Row = 41; ScrapedAt = "2026-02-03T14:49:20.053Z"; LastSeenAt = "2026-02-03T14:49:20.199Z" }
    @{ Row = 42; ScrapedAt = "2026-02-03T14:49:20.053Z"; LastSeenAt = "2026-02-03T14:49:20.201Z" }
    @{ Row = 43; ScrapedAt = "2026-02-03T14:49:20.053Z"; LastSeenAt = "2026-02-03T14:49:20.203Z" }
    @{ Row = 44; ScrapedAt = "2026-02-03T14:49:20.053Z"; LastSeenAt = "2026-02-03T14:49:20.205Z" }
    @{ Row = 45; ScrapedAt = "2026-02-03T14:49:20.053Z"; LastSeenAt = "2026-02-03T14:49:20.207Z" }
    @{ Row = 46; ScrapedAt = "2026-02-03T14:49:20.053Z"; LastSeenAt = "2026-02-03T14:49:20.215Z" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 8).Value = $u.ScrapedAt
    $ws.Cells.Item($u.Row, 9).Value = $u.LastSeenAt
}
